$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Jornada" section currently has:
#   row 7: Jornada (header)
#   row 8: Get  | Jornadas | Devuelve las jornadas
#   row 9: Post | Jornadas | Crea una jornada
#
# Insert a new "Get /Jornadas/{jornadaId}" row right after row 8, and a new
# "Put /Jornadas" row right after the (now shifted) "Post" row, pushing the
# following sections ("Partidos", "QuinelasJornada", ...) down by 2 rows.

# Insert new row 9: Get | Jornadas/{jornadaId} | Devuelve la jornada con un Id determinado
$ws.Rows("9:9").Insert()
$ws.Range("A9").Value = "Get"
$ws.Range("B9").Value = "Jornadas/{jornadaId}"
$ws.Range("C9").Value = "Devuelve la jornada con un Id determinado"

# Insert new row 11: Put | Jornadas | Abre/cierra la jornada al publico para que se creen quinelas
$ws.Rows("11:11").Insert()
$ws.Range("A11").Value = "Put"
$ws.Range("B11").Value = "Jornadas"
$ws.Range("C11").Value = "Abre/cierra la jornada al publico para que se creen quinelas"

# Restore the editor's cursor position
$ws.Range("B7").Select()
